$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ManagerList")

# Update Michael's password (row 2, column E) from "password" to "npassword"
$ws.Range("E2").Value = "npassword"
